# "Generate Report for Handback"
#
# The localization status report is regenerated: the e3329c99-...md file has
# now been handed back (status moves from "Ready for handoff" to
# "Handed back: in sync with en-US"), which also re-sorts the two
# still-pending rows (a56b4c39 / e3329c99) so a56b4c39 now sorts after
# e3329c99 on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" — File Name / zh-cn / de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md"
$ov.Range("B5").Value = "Handed back: in sync with en-US"
$ov.Range("C5").Value = "Handed back: in sync with en-US"

$ov.Range("A6").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$ov.Range("B6").Value = "Handback transform failed"
$ov.Range("C6").Value = "Handback transform failed"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c71cf16183483ba493f4f584ea26c31c72f1cf47/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/5138aaef-fea0-4fa3-addc-f35200ff812b.md", "", "", "5138aaef-fea0-4fa3-addc-f35200ff812b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f13df5d06b63f1d1c7cb021c922dd87404f41a6e/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", "a56b4c39-1622-461c-be84-e126b5128073.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 5 becomes the e3329c99 record, now fully handed back: it gains a
# "Latest Target File" and "Latest Handback File" pair and a fresh
# handoff/handback timestamp.
$zh.Range("A5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md"
$zh.Range("B5").Value = "Handed back: in sync with en-US"
$zh.Range("C5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf"
$zh.Range("D5").Value = "2016-03-09 03:24:10"
$zh.Range("E5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md"
$zh.Range("F5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf"
$zh.Range("G5").Value = "2016-03-09 03:24:56"
$zh.Range("H5").Value = "Include"
$zh.Range("E5,F5").Style = "Hyperlink"

# Row 6 becomes the a56b4c39 record (same data it always had, just moved
# down one row to make room for e3329c99 above it).
$zh.Range("A6").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$zh.Range("B6").Value = "Handback transform failed"
$zh.Range("C6").Value = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf"
$zh.Range("D6").Value = "2016-03-09 03:18:30"
$zh.Range("G6").Value = "0001-01-01 00:00:00"
$zh.Range("H6").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c71cf16183483ba493f4f584ea26c31c72f1cf47/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d8b9dfae4038afd45ca6d503c02a029c9e8709a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/437835ed84ddc4d31b435c498e3e1e30da61e6d1/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a86386962a33a2f78af6e36c6f312093de0ada21/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3d55c7d0eb607ffa99726cb296d186886cfb9af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3504c7f8f7f183c35f82b0409aa4baef8c0ec3a1/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e203a7323ea4c108dc217da020b3be848c464830/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/5138aaef-fea0-4fa3-addc-f35200ff812b.md", "", "", "5138aaef-fea0-4fa3-addc-f35200ff812b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d8b9dfae4038afd45ca6d503c02a029c9e8709a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/437835ed84ddc4d31b435c498e3e1e30da61e6d1/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a86386962a33a2f78af6e36c6f312093de0ada21/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f73095ce6d80016a9933eeb0be23a0f706063873/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6e30c0bf249e3fbb0c0f33adacb9caba30860fc5/e2e/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e30c0bf249e3fbb0c0f33adacb9caba30860fc5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f13df5d06b63f1d1c7cb021c922dd87404f41a6e/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", "a56b4c39-1622-461c-be84-e126b5128073.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8c88013cb96fb74dc99e32dd51e1385f560703b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf", "", "", "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md"
$de.Range("B5").Value = "Handed back: in sync with en-US"
$de.Range("C5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf"
$de.Range("D5").Value = "2016-03-09 03:24:21"
$de.Range("E5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md"
$de.Range("F5").Value = "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf"
$de.Range("G5").Value = "2016-03-09 03:25:30"
$de.Range("H5").Value = "Include"
$de.Range("E5,F5").Style = "Hyperlink"

$de.Range("A6").Value = "a56b4c39-1622-461c-be84-e126b5128073.md"
$de.Range("B6").Value = "Handback transform failed"
$de.Range("C6").Value = "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf"
$de.Range("D6").Value = "2016-03-09 03:18:41"
$de.Range("G6").Value = "0001-01-01 00:00:00"
$de.Range("H6").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c71cf16183483ba493f4f584ea26c31c72f1cf47/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c70354135631acf6b70af5a5bed7ba7f6dd68896/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5b773253b4b7fc5b05c8bdf26b34261fb0516739/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa3b33a2d79c20b4891a14f96f9f1d91f700fe7d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/73a7bb42a5b53ae28cff75b041d1b723ec9cc70f/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8aa48bd59e50644c8a69fd21ea7db36aca6b67a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5c22e97140677a61ceb971650cf23cd623a0eeea/e2e/36491f5a-d66c-495f-9f55-eaba4cdc0280.md", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cd3e808463306f071d0b40f4a59fbbaa5d7f7ad0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf", "", "", "36491f5a-d66c-495f-9f55-eaba4cdc0280.ff4b00ac91e37a6d18e83b1f11acee01a8980897.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/5138aaef-fea0-4fa3-addc-f35200ff812b.md", "", "", "5138aaef-fea0-4fa3-addc-f35200ff812b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c70354135631acf6b70af5a5bed7ba7f6dd68896/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5b773253b4b7fc5b05c8bdf26b34261fb0516739/e2e/257d61c9-05a1-4dd9-a061-6048d13e2c79.md", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fa3b33a2d79c20b4891a14f96f9f1d91f700fe7d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf", "", "", "257d61c9-05a1-4dd9-a061-6048d13e2c79.eb7b9669893ab9827ef85224694b2a5854a70170.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/e2e/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d12381062b24ced486cc1b013d1d816b111bfc32/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6e30c0bf249e3fbb0c0f33adacb9caba30860fc5/e2e/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e30c0bf249e3fbb0c0f33adacb9caba30860fc5/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf", "", "", "e3329c99-7ce5-47c1-8d7a-a5a0ccecc9a2.6e30c0bf249e3fbb0c0f33adacb9caba30860fc5.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f13df5d06b63f1d1c7cb021c922dd87404f41a6e/e2e/a56b4c39-1622-461c-be84-e126b5128073.md", "", "", "a56b4c39-1622-461c-be84-e126b5128073.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/924c629a31564eb8d175373e4e0f07d1ce0c4294/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf", "", "", "a56b4c39-1622-461c-be84-e126b5128073.88956322b8411ff2e30b6ae6d7edcb235bc6804c.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/d609876317ec3fcfd17bbd95305cedaff900cf52/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Report regenerated for handback of e3329c99"
